$wb = $excel.ActiveWorkbook

# --- Overview sheet: update Status for the d9b828b6... row (row 7) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B7").Value = "Handback transform failed"
$wsOverview.Range("C7").Value = "Handback transform failed"

# --- zh-cn sheet: update Status and add Error Detail for row 7 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C7").Value = "Handback transform failed"
$wsZhCn.Range("L7").Value = "Handback file name: casoa1ro.5r0 is different with handoff file name: d9b828b6-fda0-461a-9e1a-d2f6094b6f90.4aabb42fca427cb606f22ac44a4baf396b7f16e0.zh-cn."

# --- de-de sheet: update Status and add Error Detail for row 7 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C7").Value = "Handback transform failed"
$wsDeDe.Range("L7").Value = "Handback file name: casoa1ro.5r0 is different with handoff file name: d9b828b6-fda0-461a-9e1a-d2f6094b6f90.4aabb42fca427cb606f22ac44a4baf396b7f16e0.de-de."
